# "cleaned version of GA" — reshuffle the CM Freshmen timetable blocks.
#
# Strategy: each course entry is a vertically-merged block (anchor cell +
# merge range + fill color + Arial/14/bold/centered/wrap font + a medium
# box border drawn with BorderAround). Blocks that move/disappear are
# fully unmerged + cleared; blocks that are new are created from scratch;
# blocks whose anchor cell stays the same merge range just get new text
# and/or a new fill color.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- helpers -------------------------------------------------------

function Clear-Block([string]$range) {
    $r = $ws.Range($range)
    $r.UnMerge()
    $r.ClearFormats()
    $r.ClearContents()
}

function Set-Block([string]$anchor, [string]$mergeRange, [string]$text, [double]$colorBGR) {
    $ws.Range($mergeRange).Merge() | Out-Null
    $cell = $ws.Range($anchor)
    $cell.Value = $text
    $cell.Interior.Color = $colorBGR
    $cell.Font.Name = "Arial"
    $cell.Font.Bold = $true
    $cell.Font.Size = 14
    $cell.Font.Color = 0
    $cell.HorizontalAlignment = -4108   # xlCenter
    $cell.VerticalAlignment = -4108     # xlCenter
    $cell.WrapText = $true
    $ws.Range($mergeRange).BorderAround(1, -4138) | Out-Null   # xlContinuous, xlMedium
}

# BGR integer color values (Interior.Color is a BGR-packed long)
$GREEN  = 5296274    # 92D050
$GREY   = 13553360   # D0CECE
$WHITE  = 16777215   # FFFFFF
$YELLOW = 65535       # FFFF00
$RED    = 255         # FF0000

# ---- remove blocks that disappear or relocate -----------------------

Clear-Block "F6:F11"
Clear-Block "B14:B19"
Clear-Block "E14:E21"
Clear-Block "F14:F19"
Clear-Block "B26:B31"
Clear-Block "E26:E31"

# E6:E11 and C14:C19 and E34:E39 keep their merge ranges but change
# content/color below via Set-Block (Set-Block re-merges harmlessly).
# B34:B39 will be re-merged (extended) to B34:B41 below.

# ---- Monday 09:00 row -------------------------------------------------
Set-Block "B6" "B6:B11" "Mathematics II`n09:00-10:30`nroom:204" $GREEN
Set-Block "C6" "C6:C11" "Sociology`n09:00-10:30`nroom:203" $GREY
Set-Block "E6" "E6:E11" "Mathematics II`n09:00-10:30`nroom:202" $WHITE

# ---- Tuesday 11:00 row -------------------------------------------------
Set-Block "C14" "C14:C19" "Russian Language (Intermediate Level)`n11:00-12:30`nroom:204" $GREEN

# ---- Wednesday 14:00 row -----------------------------------------------
Set-Block "D26" "D26:D31" "Media Literacy`n14:00-15:30`nroom:202" $WHITE

# ---- 16:00 row -----------------------------------------------------
Set-Block "B34" "B34:B41" "Physical training`n16:00-18:00`nroom:100" $WHITE
Set-Block "C34" "C34:C39" "Russian Language (Elementary Level)`n16:00-17:30`nroom:203" $GREY
Set-Block "D34" "D34:D39" "Russian Language (Beginner Level)`n16:00-17:30`nroom:209" $YELLOW
Set-Block "E34" "E34:E39" "Media Literacy`n16:00-17:30`nroom:209" $YELLOW
Set-Block "F34" "F34:F39" "Sociology`n16:00-17:30`nroom:201" $RED
